# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column BB (06-aug) with 24 hourly values
#  - "Gaz" sheet: append a new row 51 for 2025-08-04
#  - "CO2" sheet: append a new row 51 for 2025-08-04

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" - new column BB ("06-aug")
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (BA1) onto the new header
# cell (BB1) so it keeps the bold / bordered / centered header style, then
# set its text.
$wsPrix.Range("BA1").Copy()
$wsPrix.Range("BB1").PasteSpecial(-4122)
$wsPrix.Range("BB1").Value = "06-aug"

$bbValues = @(
    75.63,
    40.07,
    34.8,
    31.1,
    41.95,
    47.55,
    51.35,
    78.89,
    83.19,
    29.71,
    -0.01,
    -0.76,
    -0.02,
    -0.02,
    -1.02,
    -0.03,
    3.18,
    52,
    80.31999999999999,
    106.22,
    114.64,
    117.02,
    110,
    90.54000000000001
)

for ($i = 0; $i -lt $bbValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 54).Value = $bbValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" - new row 51 (2025-08-04)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Plain "Value = '2025-08-04'" gets auto-recognised as a date and stored as
# a serial number, unlike the existing rows which keep the ISO date as a
# literal string. Entering it as a text-producing formula, then pasting the
# result back as a value, sidesteps Excel's date auto-detection and leaves
# the cell with the default (unstyled) format, same as the other rows.
$wsGaz.Range("A51").Formula = '="2025-08-04"'
$wsGaz.Range("A51").Copy()
$wsGaz.Range("A51").PasteSpecial(-4163)
$wsGaz.Range("B51").Value = 33.525

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" - new row 51 (2025-08-04)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A51").Formula = '="2025-08-04"'
$wsCo2.Range("A51").Copy()
$wsCo2.Range("A51").PasteSpecial(-4163)
$wsCo2.Range("B51").Value = 70.27
